$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.573.11"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.453.57"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "580.86"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").Value = "175.54"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.600"
$ws.Range("D9").Value = "3.451.63"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").Value = "6.85"
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").Value = "4.042.28"
$ws.Range("D14").Value = "30.88"
$ws.Range("E14").Value = "  -3.22%  "
$ws.Range("D15").Value = "0.132"
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").Value = "66.587.19"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "3.449.20"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D20").Value = "13.82"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "376.30"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "7.70"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "70.77"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").Value = "0.527"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  -4.60%  "
$ws.Range("D29").Value = "0.173"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "5.84"
$ws.Range("E31").Value = "  -4.89%  "
$ws.Range("D32").Value = "23.86"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("E36").Value = "  -4.47%  "
$ws.Range("E37").Value = "  -4.66%  "
$ws.Range("D38").Value = "159.62"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("D39").Value = "0.878"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "27.26"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").Value = "2.62"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").Value = "6.50"
$ws.Range("E43").Value = "  -5.19%  "
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").Value = "2.695.39"
$ws.Range("E45").Value = "  -4.41%  "
$ws.Range("D46").Value = "0.0694"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "25.24"
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("D48").Value = "40.25"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").Value = "0.0294"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "321.07"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").Value = "1.02"
$ws.Range("E51").Value = "  -2.80%  "
